$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.455.42'
$ws.Range("E2").Value = '  -2.29%  '
$ws.Range("D3").Value = '3.348.41'
$ws.Range("E3").Value = '  -4.17%  '
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").Value = '''555.35'
$ws.Range("E5").Value = '  -3.78%  '
$ws.Range("D6").Value = '''175.64'
$ws.Range("E6").Value = '  -0.85%  '
$ws.Range("D7").Value = '''0.617'
$ws.Range("E7").Value = '  -2.28%  '
$ws.Range("D8").Value = '3.340.44'
$ws.Range("E8").Value = '  -4.25%  '
$ws.Range("E9").Value = '  -0.01%  '
$ws.Range("E10").Value = '  -1.29%  '
$ws.Range("E11").Value = '  +1.16%  '
$ws.Range("D12").Value = '''54.78'
$ws.Range("E12").Value = '  -1.43%  '
$ws.Range("E13").Value = '  -2.04%  '
$ws.Range("E14").Value = '  -2.10%  '
$ws.Range("D15").Value = '3.885.04'
$ws.Range("E15").Value = '  -4.14%  '
$ws.Range("D16").Value = '''18.28'
$ws.Range("E16").Value = '  +0.15%  '
$ws.Range("E17").Value = '  -2.79%  '
$ws.Range("D18").Value = '3.345.61'
$ws.Range("E18").Value = '  -4.57%  '
$ws.Range("D19").Value = '64.365.84'
$ws.Range("E19").Value = '  -2.39%  '
$ws.Range("E20").Value = '  -1.73%  '
$ws.Range("E21").Value = '  -2.96%  '
$ws.Range("D22").Value = '''435.57'
$ws.Range("E22").Value = '  +6.27%  '
$ws.Range("E23").Value = '  +10.67%  '
$ws.Range("E24").Value = '  -3.94%  '
$ws.Range("B25").Value = 'Litecoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D25").Value = '''84.32'
$ws.Range("E25").Value = '  -0.56%  '
$ws.Range("B26").Value = 'InternetComputer(DFINITY)'
$ws.Range("C26").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D26").Value = '''13.43'
$ws.Range("E26").Value = '  +0.77%  '
$ws.Range("D27").Value = '''10.75'
$ws.Range("E27").Value = '  -2.53%  '
$ws.Range("E28").Value = '  -0.67%  '
$ws.Range("D29").Value = '''8.72'
$ws.Range("E29").Value = '  -4.44%  '
$ws.Range("D30").Value = '''29.73'
$ws.Range("E30").Value = '  -1.71%  '
$ws.Range("D31").Value = '''6.65'
$ws.Range("E31").Value = '  +0.22%  '
$ws.Range("D32").Value = '''11.47'
$ws.Range("E32").Value = '  -2.00%  '
$ws.Range("D33").Value = '''578.51'
$ws.Range("E33").Value = '  -2.31%  '
$ws.Range("D34").Value = '''0.107'
$ws.Range("E34").Value = '  -2.65%  '
$ws.Range("D35").Value = '''58.50'
$ws.Range("E35").Value = '  -3.91%  '
$ws.Range("E36").Value = '  +0.03%  '
$ws.Range("E37").Value = '  -8.06%  '
$ws.Range("D38").Value = '''3.50'
$ws.Range("E38").Value = '  -4.21%  '
$ws.Range("D39").Value = '''35.71'
$ws.Range("E39").Value = '  -2.55%  '
$ws.Range("D40").Value = '0.0₃0754'
$ws.Range("E40").Value = '  -4.81%  '
$ws.Range("E41").Value = '  -4.25%  '
$ws.Range("D42").Value = '3.102.81'
$ws.Range("E42").Value = '  -4.02%  '
$ws.Range("E43").Value = '  -0.29%  '
$ws.Range("E44").Value = '  -4.57%  '
$ws.Range("D45").Value = '''3.25'
$ws.Range("E45").Value = '  -2.31%  '
$ws.Range("D46").Value = '''0.0410'
$ws.Range("E46").Value = '  -2.03%  '
$ws.Range("E47").Value = '  -2.95%  '
$ws.Range("E48").Value = '  -1.88%  '
$ws.Range("E49").Value = '  -3.33%  '
$ws.Range("D50").Value = '''137.25'
$ws.Range("E50").Value = '  -0.32%  '
$ws.Range("E51").Value = '  -3.22%  '
